$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 5.990999999999999
$ws.Range("E2").Value = 16.914
$ws.Range("A3").Value = -21.347
$ws.Range("C3").Value = -12.09
$ws.Range("E6").Value = 16.864
$ws.Range("C12").Value = -11.536
$ws.Range("A14").Value = -21.54
$ws.Range("A16").Value = -21.368
$ws.Range("B18").Value = 5.626
$ws.Range("E19").Value = 16.298
$ws.Range("A21").Value = -20.587
$ws.Range("A23").Value = -20.686
$ws.Range("B24").Value = 5.762
$ws.Range("C24").Value = -13.134
$ws.Range("E24").Value = 16.692
$ws.Range("A25").Value = -20.761
$ws.Range("B25").Value = 6.918000000000001
$ws.Range("C25").Value = -12.577
$ws.Range("A26").Value = -21.509
$ws.Range("B27").Value = 5.546000000000001
$ws.Range("E27").Value = 16.716
$ws.Range("A29").Value = -21.104
$ws.Range("B30").Value = 6.292
$ws.Range("E30").Value = 16.139
$ws.Range("B31").Value = 5.928
$ws.Range("E31").Value = 16.4
$ws.Range("E33").Value = 17.495
$ws.Range("B39").Value = 7.894
$ws.Range("A40").Value = -20.29
$ws.Range("C41").Value = -12.61
$ws.Range("B42").Value = 7.697999999999999
$ws.Range("E42").Value = 16.555
$ws.Range("B48").Value = 5.24
$ws.Range("C50").Value = -12.962
$ws.Range("B51").Value = 5.331
$ws.Range("B52").Value = 4.972
$ws.Range("A53").Value = -21.04
$ws.Range("C53").Value = -12.269
$ws.Range("B55").Value = 4.609
$ws.Range("E55").Value = 16.499
$ws.Range("B56").Value = 5.717000000000001
$ws.Range("C56").Value = -12.728
$ws.Range("A57").Value = -21.613
$ws.Range("B57").Value = 5.672
$ws.Range("C57").Value = -13.028
$ws.Range("C58").Value = -13.081
$ws.Range("E58").Value = 16.849
$ws.Range("A59").Value = -22.273
$ws.Range("B60").Value = 5.856
$ws.Range("C61").Value = -13.181
$ws.Range("C63").Value = -11.849
$ws.Range("C64").Value = -11.897
$ws.Range("A65").Value = -21.233
$ws.Range("E65").Value = 17.407
$ws.Range("A69").Value = -21.462
$ws.Range("C70").Value = -11.969
$ws.Range("E70").Value = 17.357
$ws.Range("C72").Value = -11.828
$ws.Range("B73").Value = 6.383999999999999
$ws.Range("B74").Value = 7.575
$ws.Range("E74").Value = 16.467
$ws.Range("E75").Value = 16.936
$ws.Range("A79").Value = -21.22
$ws.Range("A83").Value = -21.14
$ws.Range("E83").Value = 16.985
$ws.Range("E84").Value = 16.438
$ws.Range("C86").Value = -13.406
$ws.Range("E86").Value = 16.583
$ws.Range("B89").Value = 5.575
$ws.Range("C89").Value = -12.225
$ws.Range("B90").Value = 5.906999999999999
$ws.Range("A91").Value = -21.53
$ws.Range("B92").Value = 6.027
$ws.Range("A93").Value = -21.324
$ws.Range("E96").Value = 16.303
$ws.Range("E97").Value = 17.228
$ws.Range("C98").Value = -12.563
$ws.Range("A100").Value = -21.291
$ws.Range("C100").Value = -13.01
$ws.Range("C102").Value = -13.237
